$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'293.72"
$ws.Range("E2").Value = "'-3.69%"
$ws.Range("G2").Value = "'8"
$ws.Range("D3").Value = "'40.56"
$ws.Range("E3").Value = "'-1.39%"
$ws.Range("G3").Value = "'8"
$ws.Range("D4").Value = "'5.019"
$ws.Range("E4").Value = "'-2.19%"
$ws.Range("G4").Value = "'8"
$ws.Range("D5").Value = "'0.07314"
$ws.Range("E5").Value = "'-3.77%"
$ws.Range("G5").Value = "'8"
$ws.Range("D6").Value = "'1.533"
$ws.Range("E6").Value = "'-8.25%"
$ws.Range("G6").Value = "'8"
$ws.Range("D7").Value = "'0.9286"
$ws.Range("E7").Value = "'-1.01%"
$ws.Range("G7").Value = "'8"
$ws.Range("D8").Value = "'2.348"
$ws.Range("E8").Value = "'-3.13%"
$ws.Range("G8").Value = "'8"
$ws.Range("D9").Value = "'0.1169"
$ws.Range("E9").Value = "'-2.75%"
$ws.Range("G9").Value = "'8"
$ws.Range("D10").Value = "'0.1752"
$ws.Range("E10").Value = "'-3.76%"
$ws.Range("G10").Value = "'8"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.04347"
$ws.Range("E11").Value = "'4.93%"
$ws.Range("G11").Value = "'8"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08724"
$ws.Range("E12").Value = "'-3.26%"
$ws.Range("G12").Value = "'8"
$ws.Range("D13").Value = "'0.1055"
$ws.Range("E13").Value = "'0.14%"
$ws.Range("G13").Value = "'8"
$ws.Range("D14").Value = "'0.001270"
$ws.Range("E14").Value = "'-1.62%"
$ws.Range("G14").Value = "'8"
$ws.Range("D15").Value = "'0.006033"
$ws.Range("E15").Value = "'3.49%"
$ws.Range("G15").Value = "'8"
$ws.Range("D16").Value = "'3.339"
$ws.Range("E16").Value = "'-0.05%"
$ws.Range("G16").Value = "'8"
$ws.Range("D17").Value = "'4.291"
$ws.Range("E17").Value = "'-0.86%"
$ws.Range("G17").Value = "'8"
$ws.Range("D18").Value = "'0.3258"
$ws.Range("E18").Value = "'-2.90%"
$ws.Range("G18").Value = "'8"
$ws.Range("D19").Value = "'7.976"
$ws.Range("E19").Value = "'4.25%"
$ws.Range("G19").Value = "'8"
$ws.Range("E20").Value = "'3.65%"
$ws.Range("G20").Value = "'8"
$ws.Range("D21").Value = "'0.2772"
$ws.Range("E21").Value = "'-2.31%"
$ws.Range("G21").Value = "'8"
$ws.Range("D22").Value = "'0.03941"
$ws.Range("E22").Value = "'0.80%"
$ws.Range("G22").Value = "'8"
$ws.Range("E23").Value = "'-1.53%"
$ws.Range("G23").Value = "'8"
$ws.Range("D24").Value = "'0.003657"
$ws.Range("E24").Value = "'-8.09%"
$ws.Range("G24").Value = "'8"
$ws.Range("E25").Value = "'-5.34%"
$ws.Range("G25").Value = "'8"
$ws.Range("E26").Value = "'22.48%"
$ws.Range("G26").Value = "'8"
$ws.Range("G27").Value = "'8"
$ws.Range("G28").Value = "'8"
$ws.Range("G29").Value = "'8"
$ws.Range("G30").Value = "'8"
$ws.Range("G31").Value = "'8"
$ws.Range("G32").Value = "'8"
$ws.Range("G33").Value = "'8"
$ws.Range("G34").Value = "'8"
$ws.Range("G35").Value = "'8"
$ws.Range("G36").Value = "'8"
$ws.Range("G37").Value = "'8"
$ws.Range("D38").Value = "'0.02316"
$ws.Range("G38").Value = "'8"
$ws.Range("D39").Value = "'0.05077"
$ws.Range("E39").Value = "'-1.23%"
$ws.Range("G39").Value = "'8"
$ws.Range("D40").Value = "'0.006407"
$ws.Range("E40").Value = "'75.32%"
$ws.Range("G40").Value = "'8"
$ws.Range("D41").Value = "'0.007856"
$ws.Range("E41").Value = "'1.75%"
$ws.Range("G41").Value = "'8"
$ws.Range("D42").Value = "'0.1289"
$ws.Range("E42").Value = "'-0.99%"
$ws.Range("G42").Value = "'8"
$ws.Range("D43").Value = "'0.007338"
$ws.Range("E43").Value = "'-3.34%"
$ws.Range("G43").Value = "'8"
$ws.Range("D44").Value = "'0.007256"
$ws.Range("E44").Value = "'-12.03%"
$ws.Range("G44").Value = "'8"
$ws.Range("D45").Value = "'0.3205"
$ws.Range("E45").Value = "'-2.58%"
$ws.Range("G45").Value = "'8"
$ws.Range("E46").Value = "'-7.06%"
$ws.Range("G46").Value = "'8"
$ws.Range("E47").Value = "'-0.14%"
$ws.Range("G47").Value = "'8"
$ws.Range("D48").Value = "'0.04715"
$ws.Range("E48").Value = "'-82.71%"
$ws.Range("G48").Value = "'8"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'-0.14%"
$ws.Range("G49").Value = "'8"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'-0.14%"
$ws.Range("G50").Value = "'8"
$ws.Range("G51").Value = "'8"
